$d = $word.ActiveDocument

# --- Step 1: locate the paragraph ending with "Результирующий массив..." ---
$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("сохраняется в файле и выводится на экран средствами ЯВУ.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "anchor paragraph not found" }
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -le $anchorRange.Start -and $pp.Range.End -ge $anchorRange.End) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq 0) { throw "could not resolve anchor paragraph index" }

# --- Step 2: replace that paragraph so it no longer carries the trailing page-break run ---
$d.Paragraphs.Item($anchorIndex).Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="ru-RU" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>Результирующий массив частотного распределения чисел по интервалам, сформированный на ассемблерном уровне, возвращается в программу, реализованную на ЯВУ, и затем сохраняется в файле и выводится на экран средствами ЯВУ.</w:t></w:r></w:p>')

# --- Step 3: insert the 8 new paragraphs (conditions block) after it ---
# new paragraph 1/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:r></w:p>')
# new paragraph 2/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="ru-RU" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>Вариант 6.</w:t></w:r></w:p>')
# new paragraph 3/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="ru-RU" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>Условия:</w:t></w:r></w:p>')
# new paragraph 4/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>- нормальное распределение;</w:t></w:r></w:p>')
$d.Paragraphs.Item($anchorIndex).Range.ParagraphFormat.FirstLineIndent = -0.01
# new paragraph 5/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>- число ассемблерных процедур — 1;</w:t></w:r></w:p>')
$d.Paragraphs.Item($anchorIndex).Range.ParagraphFormat.FirstLineIndent = -0.01
# new paragraph 6/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>- число интервалов должно быть меньше диапазона;</w:t></w:r></w:p>')
$d.Paragraphs.Item($anchorIndex).Range.ParagraphFormat.FirstLineIndent = -0.01
# new paragraph 7/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>- левые границы могут быть меньше левойграницы диапазона;</w:t></w:r></w:p>')
$d.Paragraphs.Item($anchorIndex).Range.ParagraphFormat.FirstLineIndent = -0.01
# new paragraph 8/8
$anchorP = $d.Paragraphs.Item($anchorIndex)
$anchorP.Range.InsertParagraphAfter()
$anchorIndex = $anchorIndex + 1
$newP = $d.Paragraphs.Item($anchorIndex)
$newP.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>- правая граница может быть больше правой границы диапазона.</w:t></w:r><w:r><w:br w:type="page"/></w:r></w:p>')
$d.Paragraphs.Item($anchorIndex).Range.ParagraphFormat.FirstLineIndent = -0.01

# --- Step 4: merge the tab + text runs in the "Пример работы программы." paragraph ---
$exRange = $d.Content
$foundEx = $exRange.Find.Execute("Пример работы программы.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEx) { throw "example paragraph not found" }
$exIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -le $exRange.Start -and $pp.Range.End -ge $exRange.End) {
        $exIndex = $i
        break
    }
}
if ($exIndex -eq 0) { throw "could not resolve example paragraph index" }
$d.Paragraphs.Item($exIndex).Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:jc w:val="both"/><w:rPr/></w:pPr><w:r><w:rPr><w:b/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t>Пример работы программы.</w:t></w:r></w:p>')
$d.Paragraphs.Item($exIndex).Range.ParagraphFormat.FirstLineIndent = -0.01

# --- Step 5: give the drawing-bearing run its bold/size/color rPr ---
$drawIndex = $exIndex + 1
$d.Paragraphs.Item($drawIndex).Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/><w:ind w:hanging="0"/><w:jc w:val="both"/><w:rPr><w:b/><w:b/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="0" allowOverlap="1" relativeHeight="2"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>635</wp:posOffset></wp:positionV><wp:extent cx="6120130" cy="2823845"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapSquare wrapText="largest"/><wp:docPr id="1" name="Image1" descr=""></wp:docPr><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name="Image1" descr=""></pic:cNvPr><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId2"></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="6120130" cy="2823845"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:anchor></w:drawing></w:r><w:r><w:br w:type="page"/></w:r></w:p>')
$d.Paragraphs.Item($drawIndex).Range.ParagraphFormat.FirstLineIndent = -0.01

Write-Output "All edits applied. Final paragraph count: $($d.Paragraphs.Count)"